# Applies the "confirmations/101_1.xlsx" edit:
#  - Removes the blank "Summary" divider row (old row 35), shifting the
#    trailing "Total ..." rows up by one.
#  - Relabels the per-category detail rows (New nominations / Confirmed /
#    Unconfirmed / Withdrawn / Rejected / Failed at ... ) so each one is
#    prefixed with its branch name (e.g. "     Civilian, New nominations"),
#    and renames the second "Civilian " section header to "Civilian (lists)".
#  - Fixes a few typos in the summary "Total ..." labels.
#  - Fixes B9, which held the literal text "7 9" instead of the number 79.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank "Summary" row -- everything below shifts up one row.
$ws.Rows(35).Delete()

# --- Civilian (counts) section ------------------------------------------
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Confirmed "
$ws.Range("A9").Value  = "     Civilian, Unconfirmed "
$ws.Range("B9").Value  = 79
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("A11").Value = "     Civilian, Rejected "
$ws.Range("A12").Value = "     Civilian, Failed at Aug adjournment"
$ws.Range("A13").Value = "     Civilian, Failed at sine die adjournment "

# --- Civilian (lists) section -------------------------------------------
$ws.Range("A14").Value = "Civilian (lists)"
$ws.Range("A15").Value = "     Civilian (lists), New nominations"
$ws.Range("A16").Value = "     Civilian (lists), Confirmed "
$ws.Range("A17").Value = "     Civilian (lists), Failed at sine die adjournment "

# --- Air Force section ----------------------------------------------------
$ws.Range("A19").Value = "     Air Force, New nominations"
$ws.Range("A20").Value = "     Air Force, Confirmed "
$ws.Range("A21").Value = "     Air Force, Unconfirmed "

# --- Army section ----------------------------------------------------------
$ws.Range("A23").Value = "     Army, New nominations"
$ws.Range("A24").Value = "     Army, Confirmed "
$ws.Range("A25").Value = "     Army, Unconfirmed "
$ws.Range("A26").Value = "     Army, Failed at sine die adjournment "

# --- Navy section ------------------------------------------------------
$ws.Range("A28").Value = "     Navy, New nominations"
$ws.Range("A29").Value = "     Navy, Confirmed "
$ws.Range("A30").Value = "     Navy, Unconfirmed "

# --- Marine Corps section -----------------------------------------------
$ws.Range("A32").Value = "     Marine Corps, New nominations"
$ws.Range("A33").Value = "     Marine Corps, Confirmed "
$ws.Range("A34").Value = "     Marine Corps, Unconfirmed "

# --- Summary totals (now rows 35-41 after the row 35 delete) -----------
$ws.Range("A35").Value = "Total new nominations"
$ws.Range("A38").Value = "Total withdrawn "
$ws.Range("A40").Value = "Total failed at Aug adjournment"
$ws.Range("A41").Value = "Total failed at sine die adjournment "
